$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the sheet (tab) name to reflect the new export timestamp
$ws.Name = "2023_07_05 16_35"

# Updated RSSI (D), Period (I), H (J), L (K) values, and a few Temperature (G) values
$ws.Range("D2").Value = -74
$ws.Range("I2").Value = 24597
$ws.Range("J2").Value = 40347
$ws.Range("K2").Value = 24509

$ws.Range("G3").Value = 28
$ws.Range("I3").Value = 23773
$ws.Range("J3").Value = 38686
$ws.Range("K3").Value = 23655

$ws.Range("D4").Value = -64
$ws.Range("I4").Value = 24526
$ws.Range("J4").Value = 40215
$ws.Range("K4").Value = 24420

$ws.Range("D5").Value = -62
$ws.Range("I5").Value = 24571
$ws.Range("J5").Value = 40303
$ws.Range("K5").Value = 24471

$ws.Range("D6").Value = -66
$ws.Range("G6").Value = 25
$ws.Range("I6").Value = 24472
$ws.Range("J6").Value = 40101
$ws.Range("K6").Value = 24372

$ws.Range("D7").Value = -68
$ws.Range("G7").Value = 23
$ws.Range("I7").Value = 24614
$ws.Range("J7").Value = 40382
$ws.Range("K7").Value = 24526

$ws.Range("D8").Value = -58
$ws.Range("G8").Value = 24
$ws.Range("I8").Value = 24540
$ws.Range("J8").Value = 40235
$ws.Range("K8").Value = 24446

$ws.Range("D9").Value = -66
$ws.Range("I9").Value = 24205
$ws.Range("J9").Value = 39552
$ws.Range("K9").Value = 24111

$ws.Range("D10").Value = -70
$ws.Range("G10").Value = 26
$ws.Range("I10").Value = 23863
$ws.Range("J10").Value = 38862
$ws.Range("K10").Value = 23757

$ws.Range("D11").Value = -72
$ws.Range("I11").Value = 24311
$ws.Range("J11").Value = 39772
$ws.Range("K11").Value = 24211

$wb.Save()
